$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "Nước" -> "Nước Thải"
$ws.Range("A3").Value = "Nước Thải"

# Row 4: add B4/D4 (C4 keeps "PP2")
$ws.Range("B4").Value = "abc1"
$ws.Range("D4").Value = 2222

# Row 5: add B5/D5 (C5 keeps "PP3")
$ws.Range("B5").Value = "sss"
$ws.Range("D5").Value = 22

# Row 7: add B7/D7 (C7 keeps "PPACC2")
$ws.Range("B7").Value = "kiss"
$ws.Range("D7").Value = 2

# Row 9: add B9/D9 (C9 keeps "BTPP11")
$ws.Range("B9").Value = "ds"
$ws.Range("D9").Value = 222

# Row 12: add B12/D12 (C12 keeps "Loc KKK"); B12 re-uses same text as B11 ("fde")
$ws.Range("B12").Value = "fde"
$ws.Range("D12").Value = 2222

# Update selection to A3
$ws.Range("A3").Select()
